# Apply the "Add files via upload" re-upload edit: the raw benchmark
# numbers for the "ArangoDB MMFILES (s)" row (and its derived ratio row)
# were regenerated, and the workbook was re-saved with "raw_tabular" as
# the active sheet / different cell selections remembered per-sheet.

$wb = $excel.ActiveWorkbook

$wsRaw = $wb.Worksheets.Item("raw")
$wsTab = $wb.Worksheets.Item("raw_tabular")

# ---------------------------------------------------------------------
# 1. "raw" sheet: rows 9-15 (ArangoDB MMFILES (s) block) and rows 23-29
#    (ArangoDB MMFILES (p) ratio block), column D holds the VALUE.
# ---------------------------------------------------------------------
$wsRaw.Range("D9").Value  = 8.164
$wsRaw.Range("D10").Value = 9.669
$wsRaw.Range("D11").Value = 19.658
$wsRaw.Range("D12").Value = 0.623
$wsRaw.Range("D13").Value = 1.025
$wsRaw.Range("D14").Value = 2.548
$wsRaw.Range("D15").Value = 1.582

$wsRaw.Range("D23").Value = 1.0224170319348778
$wsRaw.Range("D24").Value = 0.99629057187017001
$wsRaw.Range("D25").Value = 0.97437422552664188
$wsRaw.Range("D26").Value = 0.37083333333333335
$wsRaw.Range("D27").Value = 0.90788308237378201
$wsRaw.Range("D28").Value = 0.97924673328209078
$wsRaw.Range("D29").Value = 7.6796116504854375

# ---------------------------------------------------------------------
# 2. "raw_tabular" sheet: row 3 is "ArangoDB MMFILES (s)", row 5 is the
#    ratio row whose label moved from "ArangoDB MMFILES (p)" to
#    "ArangoDB MMFILES (s)" (and row 4's label from "ArangoDB (p)" to
#    "ArangoDB (s)").
# ---------------------------------------------------------------------
$wsTab.Range("B3").Value = 8.164
$wsTab.Range("C3").Value = 9.669
$wsTab.Range("D3").Value = 19.658
$wsTab.Range("E3").Value = 0.623
$wsTab.Range("F3").Value = 1.025
$wsTab.Range("G3").Value = 2.548
$wsTab.Range("H3").Value = 1.582

$wsTab.Range("A4").Value = "ArangoDB (s)"
$wsTab.Range("A5").Value = "ArangoDB MMFILES (s)"

$wsTab.Range("B5").Value = 1.0224170319348778
$wsTab.Range("C5").Value = 0.99629057187017001
$wsTab.Range("D5").Value = 0.97437422552664188
$wsTab.Range("E5").Value = 0.37083333333333335
$wsTab.Range("F5").Value = 0.90788308237378201
$wsTab.Range("G5").Value = 0.97924673328209078
$wsTab.Range("H5").Value = 7.6796116504854375

# ---------------------------------------------------------------------
# 3. Selection / active-sheet bookkeeping: when the file was re-saved,
#    "raw" was left with D22 selected (no longer the active tab) and
#    "raw_tabular" became the active tab with A2 selected.
# ---------------------------------------------------------------------
$wsRaw.Range("D22").Select() | Out-Null
$wsTab.Range("A2").Select() | Out-Null
